{"js": "// Update question (2) and question (3) in the project-proposal question list.\n// Q2 used to ask about listing type (1br/3br/condo/etc) - Sarah; it now asks\n// about price-point differences (500k vs 1.5m homes) - Sarah.\n// Q3 used to ask about residential vs commercial response to mortgage rates\n// - Ellen; it now asks how SFRs differ from all homes - Ellen.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst oldQ2 =\n  \"2) Does the strength of these effects differ for different types of real estate/different listings (1br vs 3br, single family homes, condos, etc) - Sarah \";\nconst oldQ3 =\n  \"3) How do different types of real estate (e.g., residential, commercial) respond to changes in mortgage rates? - Ellen \";\n\nconst newQ2 =\n  \"2) We expect the strength of this relationship to differ at different price points (500k homes vs 1.5 m homes)? - Sarah \";\nconst newQ3 =\n  \"3) How do single family rentals (SFR) differ from all homes? - Ellen \";\n\nlet q2Paragraph = null;\nlet q3Paragraph = null;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text === oldQ2) {\n    q2Paragraph = paragraphs.items[i];\n  } else if (text === oldQ3) {\n    q3Paragraph = paragraphs.items[i];\n  }\n}\n\nif (q2Paragraph) {\n  q2Paragraph.insertText(newQ2, Word.InsertLocation.replace);\n}\nif (q3Paragraph) {\n  q3Paragraph.insertText(newQ3, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update question (2) and question (3) in the project-proposal question list.\n# Q2 used to ask about listing type (1br/3br/condo/etc) - Sarah; it now asks\n# about price-point differences (500k vs 1.5m homes) - Sarah.\n# Q3 used to ask about residential vs commercial response to mortgage rates\n# - Ellen; it now asks how SFRs differ from all homes - Ellen.\n\n$d = $word.ActiveDocument\n\n$oldQ2 = \"2) Does the strength of these effects differ for different types of real estate/different listings (1br vs 3br, single family homes, condos, etc) - Sarah \"\n$oldQ3 = \"3) How do different types of real estate (e.g., residential, commercial) respond to changes in mortgage rates? - Ellen \"\n\n$newQ2 = \"2) We expect the strength of this relationship to differ at different price points (500k homes vs 1.5 m homes)? - Sarah \"\n$newQ3 = \"3) How do single family rentals (SFR) differ from all homes? - Ellen \"\n\nforeach ($p in $d.Paragraphs) {\n    $r = $p.Range\n    $text = $r.Text.TrimEnd(\"`r\")\n    if ($text -eq $oldQ2) {\n        $r.Find.Execute($oldQ2, $false, $false, $false, $false, $false, $true, 1, $false, $newQ2, 2) | Out-Null\n    }\n    elseif ($text -eq $oldQ3) {\n        $r.Find.Execute($oldQ3, $false, $false, $false, $false, $false, $true, 1, $false, $newQ3, 2) | Out-Null\n    }\n}\n"}
